# Update for BE data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BE")

# --- Row 9 (wave 1, panel B) is already there as row 10; the shared
# --- formula in I10 becomes a literal value and the row loses its
# --- survey_version (B) / spss_name (H) values.
$ws.Range("I10").Value = "be_wk09_20201120_pB_wv01"
$ws.Range("B10").ClearContents()
$ws.Range("H10").ClearContents()

# --- Three new BE (Panel B) rows get appended under the existing data.
$ws.Range("A11").Value = "be"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = "B"
$ws.Range("F11").Value = 2
$ws.Range("I11").Value = "be_wk10_19000100_pB_wv02"

$ws.Range("A12").Value = "be"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = "B"
$ws.Range("F12").Value = 3
$ws.Range("I12").Value = "be_wk11_19000100_pB_wv03"

$ws.Range("A13").Value = "be"
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 12
$ws.Range("E13").Value = "B"
$ws.Range("F13").Value = 4
$ws.Range("G13").Value = 43841
$ws.Range("H13").Value = "20_060765_BE2_Wave4_Final_v1_110121_IntClientUse"
$ws.Range("I13").Value = "be_wk12_20200111_pB_wv04"

# Match date formatting used by the rest of the column (copy format only).
$ws.Range("G2").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("A2:I2").Copy()
$ws.Range("A10:I13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply literal values/format after the format-only paste (PasteSpecial
# formats does not touch values, but G10/G13 need the date serials back in
# case paste disturbed anything).
$ws.Range("G10").Value = 44155
$ws.Range("G13").Value = 43841

# Fix back the values that ClearContents / paste may have reset.
$ws.Range("I10").Value = "be_wk09_20201120_pB_wv01"

# --- New temp sheet note ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("A1").Value = "Temp for Belguim"

# --- Selection / active tab: BE is now the active sheet ---
$ws.Activate()
$ws.Range("B9").Select()
